$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$r = $ws.Range("B2")
$r.Interior.PatternColorIndex = -4142
Write-Host "done"
